# Design_Calculations.xlsx edit:
#  - Insert a new "Shutdown Latch" worksheet between "INA225" and "Discharge"
#    with a resistor-divider / comparator-threshold calc block.
#  - Misc review fixes (kept as the diff describes: no functional change to
#    the existing three sheets beyond what Excel normally re-stamps on save).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new sheet right after "INA225" (so tab order becomes
#    INA225, Shutdown Latch, Discharge, Current Budgeting).
# ---------------------------------------------------------------------------
$ina225 = $wb.Worksheets.Item("INA225")
$latch  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ina225)
$latch.Name = "Shutdown Latch"

# ---------------------------------------------------------------------------
# 2. Type the row labels first, in the same order an author filling the
#    sheet top-to-bottom would: the four header labels, then the two
#    output labels lower down, then the units, and finally "Rp" once the
#    parallel-resistance row is worked out. This keeps new shared-string
#    insertion order lined up with the source sheet.
# ---------------------------------------------------------------------------
$latch.Range("A1").Value = "Ra"
$latch.Range("A2").Value = "Rb"
$latch.Range("A3").Value = "Rc"
$latch.Range("A4").Value = "VI"
$latch.Range("A7").Value = "Vo1"
$latch.Range("A8").Value = "Vo2"
$latch.Range("C1").Value = "ohms"
$latch.Range("C4").Value = "volts"
$latch.Range("A6").Value = "Rp"

# ---------------------------------------------------------------------------
# 3. Left-hand block (columns A:C) values / formulas - first design point.
# ---------------------------------------------------------------------------
$latch.Range("B1").Value = 100000
$latch.Range("B2").Value = 15000
$latch.Range("B3").Value = 3000
$latch.Range("B4").Value = 24

$latch.Range("C2").Value = "ohms"
$latch.Range("C3").Value = "ohms"

$latch.Range("B6").Formula = "=1/((1/B2)+(1/B3))"
$latch.Range("C6").Value = "ohms"

$latch.Range("B7").Formula = "=B4*(B2/(B1+B2))"
$latch.Range("C7").Value = "volts"

$latch.Range("B8").Formula = "=B4*(B6/(B6+B1))"
$latch.Range("C8").Value = "volts"

# Separator row - shaded.
$sep1 = $latch.Range("A5:C5")
$sep1.Interior.Color = 12566463

# ---------------------------------------------------------------------------
# 4. Right-hand block (columns E:G) - second design point, reusing the
#    same row labels/units already in the shared-string table.
# ---------------------------------------------------------------------------
$latch.Range("E1").Value = "Ra"
$latch.Range("F1").Value = 56000
$latch.Range("G1").Value = "ohms"

$latch.Range("E2").Value = "Rb"
$latch.Range("F2").Value = 10000
$latch.Range("G2").Value = "ohms"

$latch.Range("E3").Value = "Rc"
$latch.Range("F3").Value = 1000
$latch.Range("G3").Value = "ohms"

$latch.Range("E4").Value = "VI"
$latch.Range("F4").Value = 24
$latch.Range("G4").Value = "volts"

$sep2 = $latch.Range("E5:G5")
$sep2.Interior.Color = 12566463

$latch.Range("E6").Value = "Rp"
$latch.Range("F6").Formula = "=1/((1/F2)+(1/F3))"
$latch.Range("G6").Value = "ohms"

$latch.Range("E7").Value = "Vo1"
$latch.Range("F7").Formula = "=F4*(F2/(F1+F2))"
$latch.Range("G7").Value = "volts"

$latch.Range("E8").Value = "Vo2"
$latch.Range("F8").Formula = "=F4*(F6/(F6+F1))"
$latch.Range("G8").Value = "volts"

# ---------------------------------------------------------------------------
# 5. Scratch-work block lower on the sheet (columns I:J).
# ---------------------------------------------------------------------------
$latch.Range("I29").Formula = "=5-1.9"
$latch.Range("J29").Formula = "=5-3"

$latch.Range("I30").Formula = "=I29/0.008"
$latch.Range("J30").Formula = "=J29/0.008"

# ---------------------------------------------------------------------------
# 6. Cosmetics on the new sheet: column B width, active cell/selection, and
#    make this the active tab (matches the saved workbook view).
# ---------------------------------------------------------------------------
$latch.Columns.Item(2).ColumnWidth = 11.166666666666666

$latch.Range("J29").Select() | Out-Null
$latch.Activate()

$wb.Save()
